$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Comments" note for the LTC3260 row (row 9, column D)
$ws.Range("D9").Value = "mounted so that pin 1 is next to the dot"

# Widen column D (Comments) so the new, longer note is readable
$ws.Columns.Item(4).ColumnWidth = 59.67

# Reflect where the user last clicked while reviewing the BOM
$ws.Range("D12").Select()
